$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-driving content: FuelGroup/Year/metric rows 2-37 (rows 32-37 are newly added)

# Row 2: Hydrogen (2030)
$ws.Cells.Item(2, 1).Value = "Hydrogen"
$ws.Cells.Item(2, 2).Value = 2030
$ws.Cells.Item(2, 3).Value = ""
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = [double]"0.0003866335395407821"
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = [double]"1.258167754132802e-09"
$ws.Cells.Item(2, 9).Value = [double]"0.0002302426784932863"
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = ""

# Row 3: Methanol (2030)
$ws.Cells.Item(3, 1).Value = "Methanol"
$ws.Cells.Item(3, 2).Value = 2030
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = ""
$ws.Cells.Item(3, 6).Value = ""
$ws.Cells.Item(3, 7).Value = ""
$ws.Cells.Item(3, 8).Value = ""
$ws.Cells.Item(3, 9).Value = ""
$ws.Cells.Item(3, 10).Value = ""
$ws.Cells.Item(3, 11).Value = ""

# Row 4: Ammonia (2030)
$ws.Cells.Item(4, 1).Value = "Ammonia"
$ws.Cells.Item(4, 2).Value = 2030
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = ""
$ws.Cells.Item(4, 7).Value = ""
$ws.Cells.Item(4, 8).Value = ""
$ws.Cells.Item(4, 9).Value = ""
$ws.Cells.Item(4, 10).Value = ""
$ws.Cells.Item(4, 11).Value = ""

# Row 5: Synthetic Gases (2030)
$ws.Cells.Item(5, 1).Value = "Synthetic Gases"
$ws.Cells.Item(5, 2).Value = 2030
$ws.Cells.Item(5, 3).Value = ""
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = ""
$ws.Cells.Item(5, 6).Value = ""
$ws.Cells.Item(5, 7).Value = ""
$ws.Cells.Item(5, 8).Value = ""
$ws.Cells.Item(5, 9).Value = ""
$ws.Cells.Item(5, 10).Value = ""
$ws.Cells.Item(5, 11).Value = ""

# Row 6: Biogenic Gases (2030)
$ws.Cells.Item(6, 1).Value = "Biogenic Gases"
$ws.Cells.Item(6, 2).Value = 2030
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(6, 4).Value = ""
$ws.Cells.Item(6, 5).Value = [double]"0.0001276004121486529"
$ws.Cells.Item(6, 6).Value = [double]"0.000129677720281579"
$ws.Cells.Item(6, 7).Value = ""
$ws.Cells.Item(6, 8).Value = ""
$ws.Cells.Item(6, 9).Value = [double]"3.452513043873888e-05"
$ws.Cells.Item(6, 10).Value = ""
$ws.Cells.Item(6, 11).Value = ""

# Row 7: Fossil Gases (2030)
$ws.Cells.Item(7, 1).Value = "Fossil Gases"
$ws.Cells.Item(7, 2).Value = 2030
$ws.Cells.Item(7, 3).Value = ""
$ws.Cells.Item(7, 4).Value = ""
$ws.Cells.Item(7, 5).Value = ""
$ws.Cells.Item(7, 6).Value = [double]"0.0014890748389714"
$ws.Cells.Item(7, 7).Value = ""
$ws.Cells.Item(7, 8).Value = ""
$ws.Cells.Item(7, 9).Value = [double]"0.0001599571264707139"
$ws.Cells.Item(7, 10).Value = ""
$ws.Cells.Item(7, 11).Value = ""

# Row 8: Synthetic Liquids (2030)
$ws.Cells.Item(8, 1).Value = "Synthetic Liquids"
$ws.Cells.Item(8, 2).Value = 2030
$ws.Cells.Item(8, 3).Value = ""
$ws.Cells.Item(8, 4).Value = ""
$ws.Cells.Item(8, 5).Value = ""
$ws.Cells.Item(8, 6).Value = ""
$ws.Cells.Item(8, 7).Value = ""
$ws.Cells.Item(8, 8).Value = ""
$ws.Cells.Item(8, 9).Value = ""
$ws.Cells.Item(8, 10).Value = ""
$ws.Cells.Item(8, 11).Value = ""

# Row 9: Biogenic Liquids (2030)
$ws.Cells.Item(9, 1).Value = "Biogenic Liquids"
$ws.Cells.Item(9, 2).Value = 2030
$ws.Cells.Item(9, 3).Value = ""
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 5).Value = ""
$ws.Cells.Item(9, 6).Value = [double]"0.007128437630251732"
$ws.Cells.Item(9, 7).Value = [double]"1.791069058325931e-05"
$ws.Cells.Item(9, 8).Value = [double]"0.003425970247539691"
$ws.Cells.Item(9, 9).Value = [double]"0.007656994072575"
$ws.Cells.Item(9, 10).Value = [double]"1.522257158030819e-06"
$ws.Cells.Item(9, 11).Value = [double]"0.003341062137162474"

# Row 10: Fossil Liquids (2030)
$ws.Cells.Item(10, 1).Value = "Fossil Liquids"
$ws.Cells.Item(10, 2).Value = 2030
$ws.Cells.Item(10, 3).Value = ""
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = ""
$ws.Cells.Item(10, 6).Value = [double]"0.0689140591071683"
$ws.Cells.Item(10, 7).Value = [double]"0.0001261397415553"
$ws.Cells.Item(10, 8).Value = [double]"0.03130859249984613"
$ws.Cells.Item(10, 9).Value = [double]"0.0477467910813372"
$ws.Cells.Item(10, 10).Value = [double]"9.214225780114264e-06"
$ws.Cells.Item(10, 11).Value = [double]"0.0326787456338598"

# Row 11: Biomass [Solid] (2030)
$ws.Cells.Item(11, 1).Value = "Biomass [Solid]"
$ws.Cells.Item(11, 2).Value = 2030
$ws.Cells.Item(11, 3).Value = ""
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = [double]"0.001701972684689266"
$ws.Cells.Item(11, 6).Value = ""
$ws.Cells.Item(11, 7).Value = ""
$ws.Cells.Item(11, 8).Value = ""
$ws.Cells.Item(11, 9).Value = ""
$ws.Cells.Item(11, 10).Value = ""
$ws.Cells.Item(11, 11).Value = ""

# Row 12: Renewable Energy Carrier (2030)
$ws.Cells.Item(12, 1).Value = "Renewable Energy Carrier"
$ws.Cells.Item(12, 2).Value = 2030
$ws.Cells.Item(12, 3).Value = ""
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = [double]"0.001242591300436453"
$ws.Cells.Item(12, 6).Value = ""
$ws.Cells.Item(12, 7).Value = ""
$ws.Cells.Item(12, 8).Value = ""
$ws.Cells.Item(12, 9).Value = ""
$ws.Cells.Item(12, 10).Value = ""
$ws.Cells.Item(12, 11).Value = ""

# Row 13: Overall Demand (2030)
$ws.Cells.Item(13, 1).Value = "Overall Demand"
$ws.Cells.Item(13, 2).Value = 2030
$ws.Cells.Item(13, 3).Value = ""
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = [double]"0.003072164397274372"
$ws.Cells.Item(13, 6).Value = [double]"0.07804788283621379"
$ws.Cells.Item(13, 7).Value = [double]"0.0001440504321385593"
$ws.Cells.Item(13, 8).Value = [double]"0.03473456400555358"
$ws.Cells.Item(13, 9).Value = [double]"0.05582851008931494"
$ws.Cells.Item(13, 10).Value = [double]"1.073648293814508e-05"
$ws.Cells.Item(13, 11).Value = [double]"0.03601980777102227"

# Row 14: Hydrogen (2040)
$ws.Cells.Item(14, 1).Value = "Hydrogen"
$ws.Cells.Item(14, 2).Value = 2040
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = ""
$ws.Cells.Item(14, 6).Value = [double]"0.0018505625988016"
$ws.Cells.Item(14, 7).Value = ""
$ws.Cells.Item(14, 8).Value = [double]"1.053225375815418e-07"
$ws.Cells.Item(14, 9).Value = [double]"0.0003234085448113575"
$ws.Cells.Item(14, 10).Value = ""
$ws.Cells.Item(14, 11).Value = ""

# Row 15: Methanol (2040)
$ws.Cells.Item(15, 1).Value = "Methanol"
$ws.Cells.Item(15, 2).Value = 2040
$ws.Cells.Item(15, 3).Value = ""
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = ""
$ws.Cells.Item(15, 6).Value = ""
$ws.Cells.Item(15, 7).Value = ""
$ws.Cells.Item(15, 8).Value = ""
$ws.Cells.Item(15, 9).Value = ""
$ws.Cells.Item(15, 10).Value = ""
$ws.Cells.Item(15, 11).Value = ""

# Row 16: Ammonia (2040)
$ws.Cells.Item(16, 1).Value = "Ammonia"
$ws.Cells.Item(16, 2).Value = 2040
$ws.Cells.Item(16, 3).Value = ""
$ws.Cells.Item(16, 4).Value = ""
$ws.Cells.Item(16, 5).Value = ""
$ws.Cells.Item(16, 6).Value = ""
$ws.Cells.Item(16, 7).Value = ""
$ws.Cells.Item(16, 8).Value = ""
$ws.Cells.Item(16, 9).Value = ""
$ws.Cells.Item(16, 10).Value = ""
$ws.Cells.Item(16, 11).Value = ""

# Row 17: Synthetic Gases (2040)
$ws.Cells.Item(17, 1).Value = "Synthetic Gases"
$ws.Cells.Item(17, 2).Value = 2040
$ws.Cells.Item(17, 3).Value = ""
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = ""
$ws.Cells.Item(17, 6).Value = [double]"6.895466904779869e-10"
$ws.Cells.Item(17, 7).Value = ""
$ws.Cells.Item(17, 8).Value = ""
$ws.Cells.Item(17, 9).Value = [double]"9.568241845163511e-11"
$ws.Cells.Item(17, 10).Value = ""
$ws.Cells.Item(17, 11).Value = ""

# Row 18: Biogenic Gases (2040)
$ws.Cells.Item(18, 1).Value = "Biogenic Gases"
$ws.Cells.Item(18, 2).Value = 2040
$ws.Cells.Item(18, 3).Value = ""
$ws.Cells.Item(18, 4).Value = ""
$ws.Cells.Item(18, 5).Value = [double]"0.0005179764388878285"
$ws.Cells.Item(18, 6).Value = [double]"0.0001565225644405462"
$ws.Cells.Item(18, 7).Value = ""
$ws.Cells.Item(18, 8).Value = ""
$ws.Cells.Item(18, 9).Value = [double]"6.574642626760885e-05"
$ws.Cells.Item(18, 10).Value = ""
$ws.Cells.Item(18, 11).Value = ""

# Row 19: Fossil Gases (2040)
$ws.Cells.Item(19, 1).Value = "Fossil Gases"
$ws.Cells.Item(19, 2).Value = 2040
$ws.Cells.Item(19, 3).Value = ""
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 6).Value = [double]"0.0008164658281451999"
$ws.Cells.Item(19, 7).Value = ""
$ws.Cells.Item(19, 8).Value = ""
$ws.Cells.Item(19, 9).Value = [double]"0.0001681108248057366"
$ws.Cells.Item(19, 10).Value = ""
$ws.Cells.Item(19, 11).Value = ""

# Row 20: Synthetic Liquids (2040)
$ws.Cells.Item(20, 1).Value = "Synthetic Liquids"
$ws.Cells.Item(20, 2).Value = 2040
$ws.Cells.Item(20, 3).Value = ""
$ws.Cells.Item(20, 4).Value = ""
$ws.Cells.Item(20, 5).Value = ""
$ws.Cells.Item(20, 6).Value = ""
$ws.Cells.Item(20, 7).Value = ""
$ws.Cells.Item(20, 8).Value = ""
$ws.Cells.Item(20, 9).Value = ""
$ws.Cells.Item(20, 10).Value = ""
$ws.Cells.Item(20, 11).Value = ""

# Row 21: Biogenic Liquids (2040)
$ws.Cells.Item(21, 1).Value = "Biogenic Liquids"
$ws.Cells.Item(21, 2).Value = 2040
$ws.Cells.Item(21, 3).Value = ""
$ws.Cells.Item(21, 4).Value = ""
$ws.Cells.Item(21, 5).Value = ""
$ws.Cells.Item(21, 6).Value = [double]"0.003167898218339361"
$ws.Cells.Item(21, 7).Value = [double]"2.919738046383545e-05"
$ws.Cells.Item(21, 8).Value = [double]"0.004171596541923978"
$ws.Cells.Item(21, 9).Value = [double]"0.0049793053757623"
$ws.Cells.Item(21, 10).Value = [double]"1.849295387162815e-06"
$ws.Cells.Item(21, 11).Value = [double]"0.003783400348455214"

# Row 22: Fossil Liquids (2040)
$ws.Cells.Item(22, 1).Value = "Fossil Liquids"
$ws.Cells.Item(22, 2).Value = 2040
$ws.Cells.Item(22, 3).Value = ""
$ws.Cells.Item(22, 4).Value = ""
$ws.Cells.Item(22, 5).Value = ""
$ws.Cells.Item(22, 6).Value = [double]"0.0189383147124399"
$ws.Cells.Item(22, 7).Value = [double]"0.0001356072351823"
$ws.Cells.Item(22, 8).Value = [double]"0.02951858216655263"
$ws.Cells.Item(22, 9).Value = [double]"0.0209796775230359"
$ws.Cells.Item(22, 10).Value = [double]"8.176801028651042e-06"
$ws.Cells.Item(22, 11).Value = [double]"0.0316918999651766"

# Row 23: Biomass [Solid] (2040)
$ws.Cells.Item(23, 1).Value = "Biomass [Solid]"
$ws.Cells.Item(23, 2).Value = 2040
$ws.Cells.Item(23, 3).Value = ""
$ws.Cells.Item(23, 4).Value = ""
$ws.Cells.Item(23, 5).Value = [double]"0.001623178834603923"
$ws.Cells.Item(23, 6).Value = ""
$ws.Cells.Item(23, 7).Value = ""
$ws.Cells.Item(23, 8).Value = ""
$ws.Cells.Item(23, 9).Value = ""
$ws.Cells.Item(23, 10).Value = ""
$ws.Cells.Item(23, 11).Value = ""

# Row 24: Renewable Energy Carrier (2040)
$ws.Cells.Item(24, 1).Value = "Renewable Energy Carrier"
$ws.Cells.Item(24, 2).Value = 2040
$ws.Cells.Item(24, 3).Value = ""
$ws.Cells.Item(24, 4).Value = ""
$ws.Cells.Item(24, 5).Value = [double]"0.005192561614440269"
$ws.Cells.Item(24, 6).Value = ""
$ws.Cells.Item(24, 7).Value = ""
$ws.Cells.Item(24, 8).Value = ""
$ws.Cells.Item(24, 9).Value = ""
$ws.Cells.Item(24, 10).Value = ""
$ws.Cells.Item(24, 11).Value = ""

# Row 25: Overall Demand (2040)
$ws.Cells.Item(25, 1).Value = "Overall Demand"
$ws.Cells.Item(25, 2).Value = 2040
$ws.Cells.Item(25, 3).Value = ""
$ws.Cells.Item(25, 4).Value = ""
$ws.Cells.Item(25, 5).Value = [double]"0.00733371688793202"
$ws.Cells.Item(25, 6).Value = [double]"0.0249297646117133"
$ws.Cells.Item(25, 7).Value = [double]"0.0001648046156461354"
$ws.Cells.Item(25, 8).Value = [double]"0.03369028403101419"
$ws.Cells.Item(25, 9).Value = [double]"0.02651624879036532"
$ws.Cells.Item(25, 10).Value = [double]"1.002609641581386e-05"
$ws.Cells.Item(25, 11).Value = [double]"0.03547530031363182"

# Row 26: Hydrogen (2050)
$ws.Cells.Item(26, 1).Value = "Hydrogen"
$ws.Cells.Item(26, 2).Value = 2050
$ws.Cells.Item(26, 3).Value = ""
$ws.Cells.Item(26, 4).Value = ""
$ws.Cells.Item(26, 5).Value = ""
$ws.Cells.Item(26, 6).Value = [double]"0.0025703545990333"
$ws.Cells.Item(26, 7).Value = ""
$ws.Cells.Item(26, 8).Value = [double]"1.785132176394185e-07"
$ws.Cells.Item(26, 9).Value = [double]"0.0005129114715790451"
$ws.Cells.Item(26, 10).Value = ""
$ws.Cells.Item(26, 11).Value = ""

# Row 27: Methanol (2050)
$ws.Cells.Item(27, 1).Value = "Methanol"
$ws.Cells.Item(27, 2).Value = 2050
$ws.Cells.Item(27, 3).Value = ""
$ws.Cells.Item(27, 4).Value = ""
$ws.Cells.Item(27, 5).Value = ""
$ws.Cells.Item(27, 6).Value = ""
$ws.Cells.Item(27, 7).Value = ""
$ws.Cells.Item(27, 8).Value = ""
$ws.Cells.Item(27, 9).Value = ""
$ws.Cells.Item(27, 10).Value = ""
$ws.Cells.Item(27, 11).Value = ""

# Row 28: Ammonia (2050)
$ws.Cells.Item(28, 1).Value = "Ammonia"
$ws.Cells.Item(28, 2).Value = 2050
$ws.Cells.Item(28, 3).Value = ""
$ws.Cells.Item(28, 4).Value = ""
$ws.Cells.Item(28, 5).Value = ""
$ws.Cells.Item(28, 6).Value = ""
$ws.Cells.Item(28, 7).Value = ""
$ws.Cells.Item(28, 8).Value = ""
$ws.Cells.Item(28, 9).Value = ""
$ws.Cells.Item(28, 10).Value = ""
$ws.Cells.Item(28, 11).Value = ""

# Row 29: Synthetic Gases (2050)
$ws.Cells.Item(29, 1).Value = "Synthetic Gases"
$ws.Cells.Item(29, 2).Value = 2050
$ws.Cells.Item(29, 3).Value = ""
$ws.Cells.Item(29, 4).Value = ""
$ws.Cells.Item(29, 5).Value = ""
$ws.Cells.Item(29, 6).Value = [double]"7.475041808045828e-09"
$ws.Cells.Item(29, 7).Value = ""
$ws.Cells.Item(29, 8).Value = ""
$ws.Cells.Item(29, 9).Value = [double]"2.015152566902412e-09"
$ws.Cells.Item(29, 10).Value = ""
$ws.Cells.Item(29, 11).Value = ""

# Row 30: Biogenic Gases (2050)
$ws.Cells.Item(30, 1).Value = "Biogenic Gases"
$ws.Cells.Item(30, 2).Value = 2050
$ws.Cells.Item(30, 3).Value = ""
$ws.Cells.Item(30, 4).Value = ""
$ws.Cells.Item(30, 5).Value = [double]"0.001322684415299623"
$ws.Cells.Item(30, 6).Value = [double]"2.818187897676415e-05"
$ws.Cells.Item(30, 7).Value = ""
$ws.Cells.Item(30, 8).Value = ""
$ws.Cells.Item(30, 9).Value = [double]"1.934670915493983e-05"
$ws.Cells.Item(30, 10).Value = ""
$ws.Cells.Item(30, 11).Value = ""

# Row 31: Fossil Gases (2050)
$ws.Cells.Item(31, 1).Value = "Fossil Gases"
$ws.Cells.Item(31, 2).Value = 2050
$ws.Cells.Item(31, 3).Value = ""
$ws.Cells.Item(31, 4).Value = ""
$ws.Cells.Item(31, 5).Value = ""
$ws.Cells.Item(31, 6).Value = [double]"5.648282864018704e-05"
$ws.Cells.Item(31, 7).Value = ""
$ws.Cells.Item(31, 8).Value = ""
$ws.Cells.Item(31, 9).Value = [double]"6.513257528425315e-05"
$ws.Cells.Item(31, 10).Value = ""
$ws.Cells.Item(31, 11).Value = ""

# Row 32: Synthetic Liquids (2050)
$ws.Cells.Item(32, 1).Value = "Synthetic Liquids"
$ws.Cells.Item(32, 2).Value = 2050
$ws.Cells.Item(32, 3).Value = ""
$ws.Cells.Item(32, 4).Value = ""
$ws.Cells.Item(32, 5).Value = ""
$ws.Cells.Item(32, 6).Value = [double]"1.249330813810234e-11"
$ws.Cells.Item(32, 7).Value = [double]"1.075181711816921e-12"
$ws.Cells.Item(32, 8).Value = [double]"1.606411534532855e-10"
$ws.Cells.Item(32, 9).Value = [double]"8.567529130007003e-11"
$ws.Cells.Item(32, 10).Value = [double]"1.266550178375058e-14"
$ws.Cells.Item(32, 11).Value = [double]"2.612182969878195e-10"

# Row 33: Biogenic Liquids (2050)
$ws.Cells.Item(33, 1).Value = "Biogenic Liquids"
$ws.Cells.Item(33, 2).Value = 2050
$ws.Cells.Item(33, 3).Value = ""
$ws.Cells.Item(33, 4).Value = ""
$ws.Cells.Item(33, 5).Value = ""
$ws.Cells.Item(33, 6).Value = [double]"0.0003073122161435503"
$ws.Cells.Item(33, 7).Value = [double]"5.217059266314125e-05"
$ws.Cells.Item(33, 8).Value = [double]"0.005465103009091464"
$ws.Cells.Item(33, 9).Value = [double]"0.0012772935523164"
$ws.Cells.Item(33, 10).Value = [double]"2.378486713250105e-06"
$ws.Cells.Item(33, 11).Value = [double]"0.005388433374178804"

# Row 34: Fossil Liquids (2050)
$ws.Cells.Item(34, 1).Value = "Fossil Liquids"
$ws.Cells.Item(34, 2).Value = 2050
$ws.Cells.Item(34, 3).Value = ""
$ws.Cells.Item(34, 4).Value = ""
$ws.Cells.Item(34, 5).Value = ""
$ws.Cells.Item(34, 6).Value = [double]"0.0008654151399991"
$ws.Cells.Item(34, 7).Value = [double]"0.0001223119831201"
$ws.Cells.Item(34, 8).Value = [double]"0.02676532891359928"
$ws.Cells.Item(34, 9).Value = [double]"0.0037821014012264"
$ws.Cells.Item(34, 10).Value = [double]"7.038188824556429e-06"
$ws.Cells.Item(34, 11).Value = [double]"0.02953301039053015"

# Row 35: Biomass [Solid] (2050)
$ws.Cells.Item(35, 1).Value = "Biomass [Solid]"
$ws.Cells.Item(35, 2).Value = 2050
$ws.Cells.Item(35, 3).Value = ""
$ws.Cells.Item(35, 4).Value = ""
$ws.Cells.Item(35, 5).Value = [double]"0.001521101885649291"
$ws.Cells.Item(35, 6).Value = ""
$ws.Cells.Item(35, 7).Value = ""
$ws.Cells.Item(35, 8).Value = ""
$ws.Cells.Item(35, 9).Value = ""
$ws.Cells.Item(35, 10).Value = ""
$ws.Cells.Item(35, 11).Value = ""

# Row 36: Renewable Energy Carrier (2050)
$ws.Cells.Item(36, 1).Value = "Renewable Energy Carrier"
$ws.Cells.Item(36, 2).Value = 2050
$ws.Cells.Item(36, 3).Value = ""
$ws.Cells.Item(36, 4).Value = ""
$ws.Cells.Item(36, 5).Value = [double]"0.01361484606639825"
$ws.Cells.Item(36, 6).Value = ""
$ws.Cells.Item(36, 7).Value = ""
$ws.Cells.Item(36, 8).Value = ""
$ws.Cells.Item(36, 9).Value = ""
$ws.Cells.Item(36, 10).Value = ""
$ws.Cells.Item(36, 11).Value = ""

# Row 37: Overall Demand (2050)
$ws.Cells.Item(37, 1).Value = "Overall Demand"
$ws.Cells.Item(37, 2).Value = 2050
$ws.Cells.Item(37, 3).Value = ""
$ws.Cells.Item(37, 4).Value = ""
$ws.Cells.Item(37, 5).Value = [double]"0.01645863236734716"
$ws.Cells.Item(37, 6).Value = [double]"0.003827754150328017"
$ws.Cells.Item(37, 7).Value = [double]"0.000174482576858423"
$ws.Cells.Item(37, 8).Value = [double]"0.03223061059654954"
$ws.Cells.Item(37, 9).Value = [double]"0.005656787810388897"
$ws.Cells.Item(37, 10).Value = [double]"9.416675550472036e-06"
$ws.Cells.Item(37, 11).Value = [double]"0.03492144402592725"
